# Rename the usage date / count header fields for consistency:
#   R4_Month -> Usage_Date
#   R4_Count -> Usage_Count
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K1").Value = "Usage_Date"
$ws.Range("L1").Value = "Usage_Count"

# Update the selection to match the edited header cells
$ws.Range("K1:L1").Select()
